$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("G1").Value = "can_be_sold"
$ws.Range("I1").Value = "inventory_status"

# Delete the dummy data rows (rows 2 and 3), leaving only the header row
$ws.Range("A2:I3").EntireRow.Delete()
